$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUESTOS")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Add the new "ACTUALIZ" header column (H1)
$ws.Range("H1").Value = "ACTUALIZ"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Set column H width (closest achievable to the target 10.140625)
$ws.Columns.Item(8).ColumnWidth = 9.3

# Add list validation (SI,NO) for the new column, matching E/F columns
$ws.Range("H2:H1583").Validation.Add(3, 1, 1, '"SI,NO"')

# Update the active selection to match the new cursor location
[void]$ws.Range("K3").Select()
